$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C2 holds the execution time for "js in browser" (jsbr).
# Update it from "3h 23m" to "3h 35m" per the commit: "jsbr ET 3h 35m"
$ws.Range("C2").Value = "3h 35m"
